$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
$ws.Columns.Item(4).ColumnWidth = 16.8
$ws.Columns.Item(6).ColumnWidth = 7.2

# --- Title / label text changes ---
$ws.Range("A2").Value = "LAPORAN KOMISI GRO BULAN 6 TAHUN 2025"
$ws.Range("A5").Value = "Gro :"

# --- Insert 6 new rows before the TOTAL row (old row 9), pushing TOTAL to 15 and
#     TOTAL SELURUH KOMISI to 18 ---
$ws.Rows.Item(9).Resize(6).Insert()

# Copy formatting from the existing data row (row 8) down into the 6 new rows
# so the new rows 9-14 match the styles used by data rows (s=6/7/8).
$ws.Range("A8:H8").Copy()
$ws.Range("A9:H14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update existing data rows 6-8 ---
$ws.Range("A6").Value = "TF0068"
$ws.Range("B6").Value = "15-06-2025 14:33:27"
$ws.Range("C6").Value = "M019"
$ws.Range("D6").Value = "Paket Jamail"
$ws.Range("E6").Formula = "'350.000"
$ws.Range("F6").Value = 1
$ws.Range("G6").Formula = "'70.000"
$ws.Range("H6").Value = "Sdu"

$ws.Range("A7").Value = "TF0069"
$ws.Range("B7").Value = "15-06-2025 14:35:15"
$ws.Range("C7").Value = "M019"
$ws.Range("D7").Value = "Paket Jamail"
$ws.Range("E7").Formula = "'350.000"
$ws.Range("F7").Value = 1
$ws.Range("G7").Formula = "'70.000"
$ws.Range("H7").Value = "Sdu"

$ws.Range("A8").Value = "TF0071"
$ws.Range("B8").Value = "15-06-2025 14:46:23"
$ws.Range("C8").Value = "M019"
$ws.Range("D8").Value = "Paket Jamail"
$ws.Range("E8").Formula = "'350.000"
$ws.Range("F8").Value = 1
$ws.Range("G8").Formula = "'70.000"
$ws.Range("H8").Value = "Sdu"

# --- New data rows 9-14 ---
$ws.Range("A9").Value = "TF0073"
$ws.Range("B9").Value = "15-06-2025 14:52:39"
$ws.Range("C9").Value = "M019"
$ws.Range("D9").Value = "Paket Jamail"
$ws.Range("E9").Formula = "'350.000"
$ws.Range("F9").Value = 1
$ws.Range("G9").Formula = "'70.000"
$ws.Range("H9").Value = "Sdu"

$ws.Range("A10").Value = "TF0074"
$ws.Range("B10").Value = "15-06-2025 15:05:29"
$ws.Range("C10").Value = "M019"
$ws.Range("D10").Value = "Paket Jamail"
$ws.Range("E10").Formula = "'350.000"
$ws.Range("F10").Value = 1
$ws.Range("G10").Formula = "'70.000"
$ws.Range("H10").Value = "Sdu"

$ws.Range("A11").Value = "TF0075"
$ws.Range("B11").Value = "15-06-2025 15:09:27"
$ws.Range("C11").Value = "M019"
$ws.Range("D11").Value = "Paket Jamail"
$ws.Range("E11").Formula = "'350.000"
$ws.Range("F11").Value = 1
$ws.Range("G11").Formula = "'70.000"
$ws.Range("H11").Value = "Sdu"

$ws.Range("A12").Value = "TF0076"
$ws.Range("B12").Value = "15-06-2025 15:13:13"
$ws.Range("C12").Value = "M019"
$ws.Range("D12").Value = "Paket Jamail"
$ws.Range("E12").Formula = "'350.000"
$ws.Range("F12").Value = 1
$ws.Range("G12").Formula = "'70.000"
$ws.Range("H12").Value = "Sdu"

$ws.Range("A13").Value = "TF0077"
$ws.Range("B13").Value = "15-06-2025 15:14:34"
$ws.Range("C13").Value = "M019"
$ws.Range("D13").Value = "Paket Jamail"
$ws.Range("E13").Formula = "'350.000"
$ws.Range("F13").Value = 1
$ws.Range("G13").Formula = "'70.000"
$ws.Range("H13").Value = "Sdu"

$ws.Range("A14").Value = "TF0083"
$ws.Range("B14").Value = "15-06-2025 18:27:47"
$ws.Range("C14").Value = "M019"
$ws.Range("D14").Value = "Paket Jamail"
$ws.Range("E14").Formula = "'350.000"
$ws.Range("F14").Value = 1
$ws.Range("G14").Formula = "'70.000"
$ws.Range("H14").Value = "Sdu"

# --- TOTAL row (now row 15): clear the F total, put new total into G ---
$ws.Range("F15").Value = ""
$ws.Range("G15").Formula = "'630.000"

# --- TOTAL SELURUH KOMISI row (now row 18) ---
$ws.Range("G18").Formula = "'630.000"

# --- Fix up cell styles that got a stray quote-prefix flag from the text-forcing
#     apostrophe trick above, so the saved style indices match the originals. ---
$ws.Range("E6").Copy()
$ws.Range("E6:E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G6").Copy()
$ws.Range("G6:G14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A18").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
